# Append a new submission row (row 14) to the "JSS 3D" worksheet,
# mirroring the existing rows' layout (Timestamp, Full Name, Admission No, AI Score).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3D")

$ws.Cells.Item(14, 1).Value = "2026-02-15 13:38:25"
$ws.Cells.Item(14, 2).Value = "muhammad musa usman"
# Admission No is stored as text (like the rows above it) even though it
# looks numeric, so force text entry with a leading apostrophe.
$ws.Cells.Item(14, 3).Value = "'113"
$ws.Cells.Item(14, 4).Value = 10
